$d = $word.ActiveDocument
$shp = $d.Shapes.Item(3)
$tf = $shp.TextFrame
$rng = $tf.TextRange
$xml = $rng.WordOpenXML
Write-Output $xml
